# Update workbook with data for 2021-11-21 (commit: "Add data for 2021-11-29")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-11-21"

# Row 8 (June) - 2021 columns
$ws.Range("T8").Value = 4
$ws.Range("U8").Value = 126
$ws.Range("V8").Value = 0.0308

# Row 11 (September) - 2021 columns
$ws.Range("T11").Value = 7
$ws.Range("U11").Value = 171
$ws.Range("V11").Value = 0.0393

# Row 13 (November) - label and data
$ws.Range("A13").Value = "November (through 11-21)"

$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 44
$ws.Range("G13").Value = 0.12

$ws.Range("I13").Value = 83
$ws.Range("J13").Value = 0.0235

$ws.Range("O13").Value = 27
$ws.Range("P13").Value = 0.1562

$ws.Range("Q13").Value = 7
$ws.Range("R13").Value = 136
$ws.Range("S13").Value = 0.049

$ws.Range("T13").Value = 2
$ws.Range("U13").Value = 147
$ws.Range("V13").Value = 0.0134

# Row 14 (Total) - data
$ws.Range("E14").Value = 58
$ws.Range("F14").Value = 478
$ws.Range("G14").Value = 0.1082

$ws.Range("I14").Value = 732
$ws.Range("J14").Value = 0.0792

$ws.Range("O14").Value = 461
$ws.Range("P14").Value = 0.1031

$ws.Range("Q14").Value = 61
$ws.Range("R14").Value = 1139
$ws.Range("S14").Value = 0.0508

$ws.Range("T14").Value = 92
$ws.Range("U14").Value = 1499
$ws.Range("V14").Value = 0.0578
